# Append a new log row (row 38) to the "Logs" sheet, mirroring the
# existing "mercado cerrado" rows already present above it:
#   A -> Dia Actualización Valor Dólar
#   B -> EstadoMercado
#   C -> Hora Consulta
#   D -> Dia Actualización Valor Dólar (value)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A38").Value = "30/12"
$ws.Range("B38").Value = "Mercado cerrado"
$ws.Range("C38").Value = "01/01/2023 02:59"
$ws.Range("D38").Value = "848,25"
